# Commit: "Update documentation to new directions"
#
# Renames the compass-direction worksheet names to the new
# up/down/left/right-ish naming scheme, and leaves the workbook with
# README as the selected/active sheet (instead of "Small multiples").

$wb = $excel.ActiveWorkbook

# 1. Rename the direction-named sheets to their new names.
$renames = @(
    @("NNW WNW", "up-left left-up"),
    @("NNE WSW", "up-right left-down"),
    @("SSE ESE", "right-down down-right"),
    @("SSW ENE", "right-up down-left"),
    @("ABOVE LEFT border", "up-ish left-ish border"),
    @("BELOW RIGHT border", "right-ish down-ish border"),
    @("ABOVE LEFT", "up-ish left-ish"),
    @("BELOW RIGHT", "right-ish down-ish")
)

foreach ($pair in $renames) {
    $oldName = $pair[0]
    $newName = $pair[1]
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $newName
}

# 2. Reset the lingering cell selections left on several sheets back to
#    their top-left cell, then finish with "README" as the active sheet
#    (previously "Small multiples" was left active/selected).
$sheetsToReset = @(
    "Tidy",
    "Pivot",
    "up-left left-up",
    "up-ish left-ish border",
    "right-ish down-ish border",
    "Small multiples"
)

foreach ($name in $sheetsToReset) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1").Select()
}

$readme = $wb.Worksheets.Item("README")
$readme.Range("A1").Select()
